$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 54-97 were missing the "plateID" (column D) and "well_pos" (column E)
# values present for rows 2-53. This adds them, matching the rest of the sheet.

$ws.Range("D54").Value = 6
$ws.Range("E54").Value = "F8"
$ws.Range("D55").Value = 5
$ws.Range("E55").Value = "E4"
$ws.Range("D56").Value = 4
$ws.Range("E56").Value = "E9"
$ws.Range("D57").Value = 1
$ws.Range("E57").Value = "H3"
$ws.Range("D58").Value = 6
$ws.Range("E58").Value = "B11"
$ws.Range("D59").Value = 2
$ws.Range("E59").Value = "H5"
$ws.Range("D60").Value = 6
$ws.Range("E60").Value = "C1"
$ws.Range("D61").Value = 4
$ws.Range("E61").Value = "D4"
$ws.Range("D62").Value = 6
$ws.Range("E62").Value = "C7"
$ws.Range("D63").Value = 5
$ws.Range("E63").Value = "F11"
$ws.Range("D64").Value = 5
$ws.Range("E64").Value = "F12"
$ws.Range("D65").Value = 6
$ws.Range("E65").Value = "D3"
$ws.Range("D66").Value = 6
$ws.Range("E66").Value = "E1"
$ws.Range("D67").Value = 6
$ws.Range("E67").Value = "F1"
$ws.Range("D68").Value = 6
$ws.Range("E68").Value = "E9"
$ws.Range("D69").Value = 1
$ws.Range("E69").Value = "A1"
$ws.Range("D70").Value = 6
$ws.Range("E70").Value = "E12"
$ws.Range("D71").Value = 4
$ws.Range("E71").Value = "F1"
$ws.Range("D72").Value = 4
$ws.Range("E72").Value = "E3"
$ws.Range("D73").Value = 6
$ws.Range("E73").Value = "C3"
$ws.Range("D74").Value = 5
$ws.Range("E74").Value = "C5"
$ws.Range("D75").Value = 5
$ws.Range("E75").Value = "H5"
$ws.Range("D76").Value = 3
$ws.Range("E76").Value = "E10"
$ws.Range("D77").Value = 1
$ws.Range("E77").Value = "D10"
$ws.Range("D78").Value = 5
$ws.Range("E78").Value = "H9"
$ws.Range("D79").Value = 5
$ws.Range("E79").Value = "H4"
$ws.Range("D80").Value = 6
$ws.Range("E80").Value = "B4"
$ws.Range("D81").Value = 4
$ws.Range("E81").Value = "G10"
$ws.Range("D82").Value = 4
$ws.Range("E82").Value = "D11"
$ws.Range("D83").Value = 4
$ws.Range("E83").Value = "E12"
$ws.Range("D84").Value = 3
$ws.Range("E84").Value = "B6"
$ws.Range("D85").Value = 2
$ws.Range("E85").Value = "E11"
$ws.Range("D86").Value = 1
$ws.Range("E86").Value = "G5"
$ws.Range("D87").Value = 1
$ws.Range("E87").Value = "E6"
$ws.Range("D88").Value = 5
$ws.Range("E88").Value = "F5"
$ws.Range("D89").Value = 4
$ws.Range("E89").Value = "H1"
$ws.Range("D90").Value = 5
$ws.Range("E90").Value = "H12"
$ws.Range("D91").Value = 2
$ws.Range("E91").Value = "E5"
$ws.Range("D92").Value = 5
$ws.Range("E92").Value = "B2"
$ws.Range("D93").Value = 2
$ws.Range("E93").Value = "F6"
$ws.Range("D94").Value = 6
$ws.Range("E94").Value = "D2"
$ws.Range("D95").Value = 2
$ws.Range("E95").Value = "E4"
$ws.Range("D96").Value = 6
$ws.Range("E96").Value = "A3"
$ws.Range("D97").Value = 6
$ws.Range("E97").Value = "G5"
